$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Match Data" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Match Data")

# Widen the four player-passport columns (A-D) from 32.83203125 to 40.83203125
$ws1.Columns.Item(1).ColumnWidth = 40.83203125
$ws1.Columns.Item(2).ColumnWidth = 40.83203125
$ws1.Columns.Item(3).ColumnWidth = 40.83203125
$ws1.Columns.Item(4).ColumnWidth = 40.83203125

# Update header row to reflect that passport codes (not names) are required
$ws1.Range("A1").Value = "Team_1_Player_1_Passport / 第一队选手一护照码"
$ws1.Range("B1").Value = "Team_1_Player_2_Passport / 第一队选手二护照码"
$ws1.Range("C1").Value = "Team_2_Player_1_Passport / 第二队选手一护照码"
$ws1.Range("D1").Value = "Team_2_Player_2_Passport / 第二队选手二护照码"

# --- Sheet 2: "Instructions" ----------------------------------------------
$ws2 = $wb.Worksheets.Item("Instructions")

# Update existing bullet points
$ws2.Range("A4").Value = "• Use player passport codes (e.g., HVGN0BW0, KGLE38K4) - NOT names"
$ws2.Range("A5").Value = "• Leave Team_1_Player_2_Passport and Team_2_Player_2_Passport empty for singles matches"

# Insert 5 new rows starting at row 14, pushing the "Validation will check" block down
$ws2.Rows.Item(14).Resize(5).Insert()

$ws2.Range("A14").Value = "IMPORTANT - Only Passport Codes:"
$ws2.Range("A15").Value = "• The system requires passport codes, not player names"
$ws2.Range("A16").Value = "• Each passport code is 8 characters (letters and numbers)"
$ws2.Range("A17").Value = "• Find passport codes from player profiles or admin panel"
$ws2.Range("A18").Value = ""

$ws2.Range("A19").Value = "Validation will check:"
$ws2.Range("A20").Value = "• All passport codes exist in the system"
$ws2.Range("A21").Value = "• Valid score formats"
$ws2.Range("A22").Value = "• No duplicate matches"
$ws2.Range("A23").Value = "• Proper date formatting"
